$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need to be pinned to
# text format first, otherwise Excel will auto-convert the literal into a
# numeric value (and we would lose the exact original formatting/precision).
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D12",
    "D13",
    "D14",
    "D15",
    "D16",
    "D18",
    "D19",
    "D22",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51",
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin price / volume figures scraped by the Actions job.
$ws.Range("D2").Value = '29.792.23'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").Value = '1.875.54'
$ws.Range("E3").Value = '  +1.67%  '
$ws.Range("D4").Value = '0.9981'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '242.83'
$ws.Range("E5").Value = '  -2.25%  '
$ws.Range("D6").Value = '0.9988'
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = '0.4932'
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("D8").Value = '43.79'
$ws.Range("E8").Value = '  -1.48%  '
$ws.Range("D9").Value = '0.2901'
$ws.Range("E9").Value = '  +3.07%  '
$ws.Range("D10").Value = '0.06601'
$ws.Range("E10").Value = '  +2.23%  '
$ws.Range("D11").Value = '1.877.95'
$ws.Range("E11").Value = '  +1.88%  '
$ws.Range("D12").Value = '16.90'
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").Value = '0.07163'
$ws.Range("E13").Value = '  +0.53%  '
$ws.Range("D14").Value = '0.6669'
$ws.Range("E14").Value = '  +1.00%  '
$ws.Range("D15").Value = '85.33'
$ws.Range("E15").Value = '  +0.90%  '
$ws.Range("D16").Value = '4.804'
$ws.Range("E16").Value = '  +1.46%  '
$ws.Range("D17").Value = '29.765.79'
$ws.Range("E17").Value = '  +0.05%  '
$ws.Range("D18").Value = '0.000007792'
$ws.Range("E18").Value = '  +5.64%  '
$ws.Range("D19").Value = '0.9994'
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("E20").Value = '  +2.28%  '
$ws.Range("D21").Value = '2.121.48'
$ws.Range("E21").Value = '  +2.42%  '
$ws.Range("D22").Value = '0.9999'
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("E23").Value = '  +2.95%  '
$ws.Range("D24").Value = '9.117'
$ws.Range("E24").Value = '  +2.56%  '
$ws.Range("D25").Value = '5.558'
$ws.Range("E25").Value = '  +2.25%  '
$ws.Range("D26").Value = '147.32'
$ws.Range("E26").Value = '  +2.89%  '
$ws.Range("D27").Value = '134.35'
$ws.Range("E27").Value = '  +2.42%  '
$ws.Range("D28").Value = '16.65'
$ws.Range("E28").Value = '  +0.86%  '
$ws.Range("D29").Value = '1.921'
$ws.Range("E29").Value = '  +1.03%  '
$ws.Range("E30").Value = '  -1.51%  '
$ws.Range("D31").Value = '4.179'
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("D32").Value = '0.08658'
$ws.Range("E32").Value = '  +1.19%  '
$ws.Range("D33").Value = '3.909'
$ws.Range("E33").Value = '  +2.34%  '
$ws.Range("D34").Value = '0.05054'
$ws.Range("E34").Value = '  +1.38%  '
$ws.Range("D35").Value = '0.7063'
$ws.Range("E35").Value = '  +4.33%  '
$ws.Range("D36").Value = '1.106'
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").Value = '2.670'
$ws.Range("E37").Value = '  -1.07%  '
$ws.Range("D38").Value = '2.205'
$ws.Range("E38").Value = '  -3.93%  '
$ws.Range("D39").Value = '2.679'
$ws.Range("E39").Value = '  -1.74%  '
$ws.Range("D40").Value = '0.9289'
$ws.Range("E40").Value = '  -2.77%  '
$ws.Range("D41").Value = '0.01636'
$ws.Range("E41").Value = '  +2.02%  '
$ws.Range("D42").Value = '6.049'
$ws.Range("E42").Value = '  -1.44%  '
$ws.Range("D43").Value = '0.9952'
$ws.Range("E43").Value = '  -0.27%  '
$ws.Range("D44").Value = '102.58'
$ws.Range("E44").Value = '  -0.47%  '
$ws.Range("D45").Value = '0.4164'
$ws.Range("E45").Value = '  +1.77%  '
$ws.Range("D46").Value = '7.539'
$ws.Range("E46").Value = '  +3.77%  '
$ws.Range("D47").Value = '0.1256'
$ws.Range("E47").Value = '  +2.02%  '
$ws.Range("D48").Value = '0.05694'
$ws.Range("E48").Value = '  +2.01%  '
$ws.Range("D49").Value = '32.54'
$ws.Range("E49").Value = '  +2.08%  '
$ws.Range("D50").Value = '8.209'
$ws.Range("E50").Value = '  +0.87%  '
$ws.Range("D51").Value = '1.339'
$ws.Range("E51").Value = '  +1.46%  '
